# The "Employment impact" row (row 69) is removed from the Specification
# sheet. Deleting the row shifts every row below it up by one, which also
# naturally contracts the used range from A1:I181 down to A1:I180 and
# re-numbers the vertical merges of the "top-level"/"top-level-description"
# (columns A/B) grouping cells that sit below the deleted row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(69).Delete()
